$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1690140845070423
$ws.Range("C2").Value = 0.5830985915492958
$ws.Range("J2").Value = 0.01690140845070422
$ws.Range("P2").Value = 0.1211267605633803
$ws.Range("S2").Value = 0.1098591549295775
$ws.Range("B3").Value = 0.01851851851851852
$ws.Range("C3").Value = 0.04629629629629629
$ws.Range("J3").Value = 0.02777777777777778
$ws.Range("P3").Value = 0.7222222222222222
$ws.Range("S3").Value = 0.1851851851851852
$ws.Range("J4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.7346938775510204
$ws.Range("S4").Value = 0.2448979591836735
$ws.Range("B6").Value = 0.07664233576642336
$ws.Range("D6").Value = 0.0072992700729927
$ws.Range("F6").Value = 0.06204379562043796
$ws.Range("J6").Value = 0.2846715328467153
$ws.Range("O6").Value = 0.03284671532846715
$ws.Range("Q6").Value = 0.1642335766423358
$ws.Range("R6").Value = 0.0583941605839416
$ws.Range("S6").Value = 0.3138686131386861
$ws.Range("B7").Value = 0.09230769230769231
$ws.Range("D7").Value = 0.01923076923076923
$ws.Range("F7").Value = 0.08461538461538462
$ws.Range("J7").Value = 0.07692307692307693
$ws.Range("O7").Value = 0.01538461538461539
$ws.Range("Q7").Value = 0.1961538461538462
$ws.Range("R7").Value = 0.06153846153846154
$ws.Range("S7").Value = 0.4538461538461538
$ws.Range("B8").Value = 0.1103603603603604
$ws.Range("D8").Value = 0.01351351351351351
$ws.Range("E8").Value = 0.002252252252252252
$ws.Range("F8").Value = 0.04954954954954955
$ws.Range("J8").Value = 0.1373873873873874
$ws.Range("O8").Value = 0.01126126126126126
$ws.Range("Q8").Value = 0.1756756756756757
$ws.Range("R8").Value = 0.1036036036036036
$ws.Range("S8").Value = 0.3963963963963964
$ws.Range("B9").Value = 0.0912863070539419
$ws.Range("D9").Value = 0.008298755186721992
$ws.Range("F9").Value = 0.1078838174273859
$ws.Range("J9").Value = 0.1410788381742739
$ws.Range("O9").Value = 0.008298755186721992
$ws.Range("Q9").Value = 0.1618257261410788
$ws.Range("R9").Value = 0.0912863070539419
$ws.Range("S9").Value = 0.3900414937759336
$ws.Range("B10").Value = 0.109159347553325
$ws.Range("D10").Value = 0.02258469259723965
$ws.Range("E10").Value = 0.001254705144291092
$ws.Range("F10").Value = 0.06524466750313676
$ws.Range("J10").Value = 0.1373902132998745
$ws.Range("O10").Value = 0.0150564617314931
$ws.Range("Q10").Value = 0.2277289836888331
$ws.Range("R10").Value = 0.08218318695106649
$ws.Range("S10").Value = 0.3393977415307403
$ws.Range("J11").Value = 0.07990314769975787
$ws.Range("K11").Value = 0.2082324455205811
$ws.Range("L11").Value = 0.549636803874092
$ws.Range("S11").Value = 0.01937046004842615
$ws.Range("G12").Value = 0.7533039647577092
$ws.Range("J12").Value = 0.1982378854625551
$ws.Range("K12").Value = 0.00881057268722467
$ws.Range("L12").Value = 0.013215859030837
$ws.Range("S12").Value = 0.02643171806167401
$ws.Range("F13").Value = 0.01818181818181818
$ws.Range("G13").Value = 0.6
$ws.Range("J13").Value = 0.3272727272727273
$ws.Range("S13").Value = 0.05454545454545454
$ws.Range("F15").Value = 0.03658536585365853
$ws.Range("H15").Value = 0.1422764227642276
$ws.Range("I15").Value = 0.08943089430894309
$ws.Range("J15").Value = 0.3699186991869919
$ws.Range("K15").Value = 0.05284552845528456
$ws.Range("M15").Value = 0.02032520325203252
$ws.Range("O15").Value = 0.04471544715447155
$ws.Range("S15").Value = 0.2439024390243902
$ws.Range("F16").Value = 0.04405286343612335
$ws.Range("H16").Value = 0.13215859030837
$ws.Range("I16").Value = 0.09251101321585903
$ws.Range("J16").Value = 0.4096916299559472
$ws.Range("K16").Value = 0.1145374449339207
$ws.Range("M16").Value = 0.03524229074889868
$ws.Range("O16").Value = 0.01762114537444934
$ws.Range("S16").Value = 0.1541850220264317
$ws.Range("F17").Value = 0.01043478260869565
$ws.Range("H17").Value = 0.1408695652173913
$ws.Range("I17").Value = 0.09043478260869565
$ws.Range("J17").Value = 0.4295652173913043
$ws.Range("K17").Value = 0.1060869565217391
$ws.Range("M17").Value = 0.01739130434782609
$ws.Range("N17").Value = 0.001739130434782609
$ws.Range("O17").Value = 0.08347826086956522
$ws.Range("S17").Value = 0.12
$ws.Range("F18").Value = 0.01762114537444934
$ws.Range("H18").Value = 0.1101321585903084
$ws.Range("I18").Value = 0.09691629955947137
$ws.Range("J18").Value = 0.4449339207048458
$ws.Range("K18").Value = 0.1541850220264317
$ws.Range("M18").Value = 0.004405286343612335
$ws.Range("N18").Value = 0.004405286343612335
$ws.Range("O18").Value = 0.07048458149779736
$ws.Range("S18").Value = 0.09691629955947137
$ws.Range("F19").Value = 0.01430517711171662
$ws.Range("H19").Value = 0.1900544959128065
$ws.Range("I19").Value = 0.08583106267029973
$ws.Range("J19").Value = 0.3801089918256131
$ws.Range("K19").Value = 0.1246594005449591
$ws.Range("M19").Value = 0.02247956403269755
$ws.Range("O19").Value = 0.06335149863760219
$ws.Range("S19").Value = 0.1192098092643052
